# Scheduled runner update: refresh cached market-board derived figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# profit sheets. Values come from an external price-refresh pass; the
# Leve definition columns (A-G) are untouched.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -888
$ws.Range("N7").ClearContents()

# Row 14
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -809
$ws.Range("N14").ClearContents()

# Row 28
$ws.Range("H28").Value = 423.29413
$ws.Range("I28").Value = 523
$ws.Range("K28").Value = 523
$ws.Range("M28").Value = -38

# Row 32
$ws.Range("H32").Value = 5686.3335
$ws.Range("I32").Value = 4382.5
$ws.Range("J32").Value = 6555.5557
$ws.Range("K32").Value = 4382.5
$ws.Range("L32").Value = 6555.5557
$ws.Range("M32").Value = -4056.5
$ws.Range("N32").Value = -7207.5557

# Row 100
$ws.Range("H100").Value = 2897
$ws.Range("I100").Value = 2897
$ws.Range("K100").Value = 2897
$ws.Range("M100").Value = -2356

# Row 132
$ws.Range("H132").Value = 11579.4
$ws.Range("I132").Value = 1849.25
$ws.Range("K132").Value = 5547.75
$ws.Range("M132").Value = -3017.75

# Row 138
$ws.Range("H138").Value = 5982.5713
$ws.Range("I138").Value = 6572.5
$ws.Range("J138").Value = 5196
$ws.Range("K138").Value = 19717.5
$ws.Range("L138").Value = 15588
$ws.Range("M138").Value = -14577.5
$ws.Range("N138").Value = -25868

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 8341452.5
$ws.Range("I61").Value = 13894930
$ws.Range("K61").Value = 13894930
$ws.Range("M61").Value = -13894718

# Row 63
$ws.Range("H63").Value = 1001958.4
$ws.Range("I63").Value = 2307
$ws.Range("J63").Value = 3334478.2
$ws.Range("K63").Value = 2307
$ws.Range("L63").Value = 3334478.2
$ws.Range("M63").Value = -1621
$ws.Range("N63").Value = -3335850.2

# Row 66
$ws.Range("H66").Value = 1001958.4
$ws.Range("I66").Value = 2307
$ws.Range("J66").Value = 3334478.2
$ws.Range("K66").Value = 11535
$ws.Range("L66").Value = 16672391
$ws.Range("M66").Value = -8103
$ws.Range("N66").Value = -16679255

# Row 102
$ws.Range("H102").Value = 2569.6052
$ws.Range("I102").Value = 2061.7742
$ws.Range("J102").Value = 4818.5713
$ws.Range("K102").Value = 2061.7742
$ws.Range("L102").Value = 4818.5713
$ws.Range("M102").Value = -439.7741999999998
$ws.Range("N102").Value = -8062.5713

# Row 110
$ws.Range("H110").Value = 2615
$ws.Range("I110").Value = 1872.2858
$ws.Range("J110").Value = 4925.6665
$ws.Range("K110").Value = 1872.2858
$ws.Range("L110").Value = 4925.6665
$ws.Range("M110").Value = 172.7141999999999
$ws.Range("N110").Value = -9015.666499999999

# Row 122
$ws.Range("H122").Value = 5823.3335
$ws.Range("I122").Value = 4750
$ws.Range("J122").Value = 6360
$ws.Range("K122").Value = 14250
$ws.Range("L122").Value = 19080
$ws.Range("M122").Value = -11800
$ws.Range("N122").Value = -23980

# Row 136
$ws.Range("H136").Value = 8341452.5
$ws.Range("I136").Value = 13894930
$ws.Range("K136").Value = 41684790
$ws.Range("M136").Value = -41682240

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 1500
$ws.Range("I75").Value = 1500
$ws.Range("K75").Value = 1500
$ws.Range("M75").Value = -564

# Row 78
$ws.Range("H78").Value = 1500
$ws.Range("I78").Value = 1500
$ws.Range("K78").Value = 4500
$ws.Range("M78").Value = 180

# Row 86
$ws.Range("H86").Value = 18183500
$ws.Range("I86").Value = 1849.9
$ws.Range("J86").Value = 200000000
$ws.Range("K86").Value = 1849.9
$ws.Range("L86").Value = 200000000
$ws.Range("M86").Value = -726.9000000000001
$ws.Range("N86").Value = -200002246

# Row 89
$ws.Range("H89").Value = 18183500
$ws.Range("I89").Value = 1849.9
$ws.Range("J89").Value = 200000000
$ws.Range("K89").Value = 9249.5
$ws.Range("L89").Value = 1000000000
$ws.Range("M89").Value = -3633.5
$ws.Range("N89").Value = -1000011232

# Row 99
$ws.Range("H99").Value = 5492.5
$ws.Range("I99").Value = 4312.7144
$ws.Range("J99").Value = 6243.273
$ws.Range("K99").Value = 4312.7144
$ws.Range("L99").Value = 6243.273
$ws.Range("M99").Value = -2814.7144
$ws.Range("N99").Value = -9239.273000000001

# Row 105
$ws.Range("H105").Value = 4543.2593
$ws.Range("I105").Value = 3970.238
$ws.Range("J105").Value = 6548.8335
$ws.Range("K105").Value = 3970.238
$ws.Range("L105").Value = 6548.8335
$ws.Range("M105").Value = -2223.238
$ws.Range("N105").Value = -10042.8335

# Row 107
$ws.Range("H107").Value = 5571.2856
$ws.Range("I107").Value = 5833.3335
$ws.Range("K107").Value = 5833.3335
$ws.Range("M107").Value = -3913.3335

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 97
$ws.Range("H97").Value = 30193.8
$ws.Range("J97").Value = 30242.25
$ws.Range("L97").Value = 30242.25
$ws.Range("N97").Value = -32224.25

# Row 104
$ws.Range("H104").Value = 33000
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 107
$ws.Range("H107").Value = 1538.24
$ws.Range("I107").Value = 734.8889
$ws.Range("J107").Value = 1990.125
$ws.Range("K107").Value = 734.8889
$ws.Range("L107").Value = 1990.125
$ws.Range("M107").Value = 1185.1111
$ws.Range("N107").Value = -5830.125

# Row 132
$ws.Range("H132").Value = 2314.75
$ws.Range("I132").Value = 1859.7142
$ws.Range("K132").Value = 5579.142599999999
$ws.Range("M132").Value = -3049.142599999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 177.77777
$ws.Range("I14").Value = 177.77777
$ws.Range("K14").Value = 533.33331
$ws.Range("M14").Value = -360.33331

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 445587.78
$ws.Range("I113").Value = 667381.7
$ws.Range("K113").Value = 667381.7
$ws.Range("M113").Value = -665211.7

# Row 122
$ws.Range("H122").Value = 2175.842
$ws.Range("I122").Value = 2287
$ws.Range("J122").Value = 1985.2858
$ws.Range("K122").Value = 6861
$ws.Range("L122").Value = 5955.857400000001
$ws.Range("M122").Value = -4411
$ws.Range("N122").Value = -10855.8574

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2599.9285
$ws.Range("I16").Value = 2783.25
$ws.Range("K16").Value = 2783.25
$ws.Range("M16").Value = -2613.25

# Row 46
$ws.Range("H46").Value = 2043
$ws.Range("I46").Value = 1598
$ws.Range("J46").Value = 2599.25
$ws.Range("K46").Value = 1598
$ws.Range("L46").Value = 2599.25
$ws.Range("M46").Value = -1410
$ws.Range("N46").Value = -2975.25

# Row 132
$ws.Range("H132").Value = 9350.875
$ws.Range("I132").Value = 11800.272
$ws.Range("J132").Value = 3962.2
$ws.Range("K132").Value = 35400.81600000001
$ws.Range("L132").Value = 11886.6
$ws.Range("M132").Value = -32870.81600000001
$ws.Range("N132").Value = -16946.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 102
$ws.Range("H102").Value = 70000
$ws.Range("J102").Value = 70000
$ws.Range("L102").Value = 70000
$ws.Range("N102").Value = -76490

# Row 136
$ws.Range("H136").Value = 5712.971
$ws.Range("I136").Value = 5295.3228
$ws.Range("J136").Value = 8949.75
$ws.Range("K136").Value = 15885.9684
$ws.Range("L136").Value = 26849.25
$ws.Range("M136").Value = -13335.9684
$ws.Range("N136").Value = -31949.25

